$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 125, pushing the blank separator + summary rows down by one.
$ws.Rows.Item(125).Insert()

# Fill the new row 125 with the new working-hours entry.
$ws.Cells.Item(125, 1).Value = 2014
$ws.Cells.Item(125, 2).Value = 4
$ws.Cells.Item(125, 3).Value = 14
$ws.Cells.Item(125, 4).Value = 0.3888888888888889
$ws.Cells.Item(125, 5).Value = 0.5

$ws.Range("F125").Formula = "=(E125-D125)*24*60"
$ws.Range("G125").Formula = "=F125/60"

$ws.Range("F125").Select()
